$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename "heldentyp" key to "profession" (value Soldatin stays in B5)
$ws.Range("A5").Value = "profession"

# Insert a new row 6 with "beruf" / "Schriftstellerin", pushing everything below down by one row
$ws.Rows("6").Insert()
$ws.Range("A6").Value = "beruf"
$ws.Range("B6").Value = "Schriftstellerin"

# After the insert, the old row 7 (MU) became row 8, KL/IN/CH became rows 9/10/11, etc.
# Update the attribute values that changed (KL, IN, CH: 8 -> 10)
$ws.Range("B9").Value = 10
$ws.Range("B10").Value = 10
$ws.Range("B11").Value = 10

# AsP value changes from 4 to 14 (now at row 18)
$ws.Range("B18").Value = 14

# Insert two new rows (KaP=0, INI=12) right after the AsP row (row 18),
# pushing the old "Int" row (12) and "AP" row down.
$ws.Rows("19:20").Insert()
$ws.Range("A19").Value = "KaP"
$ws.Range("B19").Value = 0
$ws.Range("A20").Value = "INI"
$ws.Range("B20").Value = 12

# The old "Int" row (value 12) is now at row 21; remove it entirely so that the
# "AP" row (row 22) takes its place as row 21.
$ws.Rows("21").Delete()

# Add new rows for geschwindigkeit / ausweichen after a blank separator row
$ws.Range("A23").Value = "geschwindigkeit"
$ws.Range("B23").Value = 8
$ws.Range("A24").Value = "ausweichen"
$ws.Range("B24").Value = 10

# Add new has_grimoire row at the bottom, after another blank separator row
$ws.Range("A26").Value = "has_grimoire"
$ws.Range("B26").Value = 0

# Column widths as set in the updated sheet (closest achievable values given
# this runtime's pixel-grid rounding of ColumnWidth -> stored xml width)
$ws.Columns("A").ColumnWidth = 27.6
$ws.Columns("B").ColumnWidth = 23

# Update the active selection to match the edited sheet
$ws.Range("C25").Select()
